$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28,8).Value = 1044.6666  # H28
$ws.Cells.Item(28,9).Value = 504.63635  # I28
$ws.Cells.Item(28,10).Value = 1501.6154  # J28
$ws.Cells.Item(28,11).Value = 504.63635  # K28
$ws.Cells.Item(28,12).Value = 1501.6154  # L28
$ws.Cells.Item(28,13).Value = -19.63634999999999  # M28
$ws.Cells.Item(28,14).Value = -2471.6154  # N28

$ws.Cells.Item(33,8).Value = 11210.2  # H33
$ws.Cells.Item(33,9).Value = 14574.6  # I33
$ws.Cells.Item(33,10).Value = 1117  # J33
$ws.Cells.Item(33,11).Value = 14574.6  # K33
$ws.Cells.Item(33,12).Value = 1117  # L33
$ws.Cells.Item(33,13).Value = -14345.6  # M33
$ws.Cells.Item(33,14).Value = -1575  # N33

$ws.Cells.Item(76,8).Value = 13639  # H76
$ws.Cells.Item(76,9).Value = 13639  # I76
$ws.Cells.Item(76,11).Value = 13639  # K76
$ws.Cells.Item(76,13).Value = -13324  # M76

$ws.Cells.Item(79,8).Value = 13639  # H79
$ws.Cells.Item(79,9).Value = 13639  # I79
$ws.Cells.Item(79,11).Value = 13639  # K79
$ws.Cells.Item(79,13).Value = -12547  # M79

$ws.Cells.Item(86,8).Value = 2457.8  # H86
$ws.Cells.Item(86,9).Value = 1949.5  # I86
$ws.Cells.Item(86,10).Value = 2796.6667  # J86
$ws.Cells.Item(86,11).Value = 1949.5  # K86
$ws.Cells.Item(86,12).Value = 2796.6667  # L86
$ws.Cells.Item(86,13).Value = -826.5  # M86
$ws.Cells.Item(86,14).Value = -5042.6667  # N86

$ws.Cells.Item(89,8).Value = 2457.8  # H89
$ws.Cells.Item(89,9).Value = 1949.5  # I89
$ws.Cells.Item(89,10).Value = 2796.6667  # J89
$ws.Cells.Item(89,11).Value = 9747.5  # K89
$ws.Cells.Item(89,12).Value = 13983.3335  # L89
$ws.Cells.Item(89,13).Value = -4131.5  # M89
$ws.Cells.Item(89,14).Value = -25215.3335  # N89

$ws.Cells.Item(118,8).Value = 291.7143  # H118
$ws.Cells.Item(118,9).Value = 291.7143  # I118
$ws.Cells.Item(118,11).Value = 875.1428999999999  # K118
$ws.Cells.Item(118,13).Value = 781.8571000000001  # M118

$ws.Cells.Item(129,8).Value = 2868.7778  # H129
$ws.Cells.Item(129,9).Value = 1650  # I129
$ws.Cells.Item(129,11).Value = 4950  # K129
$ws.Cells.Item(129,13).Value = 50  # M129

$ws.Cells.Item(132,8).Value = 4604.5635  # H132
$ws.Cells.Item(132,9).Value = 4956.64  # I132
$ws.Cells.Item(132,10).Value = 1083.8  # J132
$ws.Cells.Item(132,11).Value = 14869.92  # K132
$ws.Cells.Item(132,12).Value = 3251.4  # L132
$ws.Cells.Item(132,13).Value = -12339.92  # M132
$ws.Cells.Item(132,14).Value = -8311.4  # N132

$ws.Cells.Item(138,8).Value = 29414016  # H138
$ws.Cells.Item(138,9).Value = 1277.579  # I138
$ws.Cells.Item(138,10).Value = 66670150  # J138
$ws.Cells.Item(138,11).Value = 3832.737  # K138
$ws.Cells.Item(138,12).Value = 200010450  # L138
$ws.Cells.Item(138,13).Value = 1307.263  # M138
$ws.Cells.Item(138,14).Value = -200020730  # N138

$ws.Cells.Item(141,8).Value = 1486.4828  # H141
$ws.Cells.Item(141,9).Value = 1231.76  # I141
$ws.Cells.Item(141,10).Value = 3078.5  # J141
$ws.Cells.Item(141,11).Value = 3695.28  # K141
$ws.Cells.Item(141,12).Value = 9235.5  # L141
$ws.Cells.Item(141,13).Value = 1484.72  # M141
$ws.Cells.Item(141,14).Value = -19595.5  # N141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74,8).Value = 50911.76  # H74
$ws.Cells.Item(74,9).Value = 60395.76  # I74
$ws.Cells.Item(74,11).Value = 60395.76  # K74
$ws.Cells.Item(74,13).Value = -59521.76  # M74

$ws.Cells.Item(77,8).Value = 50911.76  # H77
$ws.Cells.Item(77,9).Value = 60395.76  # I77
$ws.Cells.Item(77,11).Value = 301978.8  # K77
$ws.Cells.Item(77,13).Value = -297610.8  # M77

$ws.Cells.Item(97,8).Value = 1816.3077  # H97
$ws.Cells.Item(97,9).Value = 1345.8889  # I97
$ws.Cells.Item(97,10).Value = 2874.75  # J97
$ws.Cells.Item(97,11).Value = 1345.8889  # K97
$ws.Cells.Item(97,12).Value = 2874.75  # L97
$ws.Cells.Item(97,13).Value = -849.8888999999999  # M97
$ws.Cells.Item(97,14).Value = -3866.75  # N97

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(44,8).Value = 25000  # H44
$ws.Cells.Item(44,9).Value = 25000  # I44
$ws.Cells.Item(44,11).Value = 25000  # K44
$ws.Cells.Item(44,13).Value = -24503  # M44

$ws.Cells.Item(86,8).Value = 14478.777  # H86
$ws.Cells.Item(86,10).Value = 20656.584  # J86
$ws.Cells.Item(86,12).Value = 20656.584  # L86
$ws.Cells.Item(86,14).Value = -22902.584  # N86

$ws.Cells.Item(89,8).Value = 14478.777  # H89
$ws.Cells.Item(89,10).Value = 20656.584  # J89
$ws.Cells.Item(89,12).Value = 103282.92  # L89
$ws.Cells.Item(89,14).Value = -114514.92  # N89

$ws.Cells.Item(105,8).Value = 50937.25  # H105
$ws.Cells.Item(105,9).Value = 50937.25  # I105
$ws.Cells.Item(105,10).Value = 0  # J105
$ws.Cells.Item(105,11).Value = 50937.25  # K105
$ws.Cells.Item(105,12).Value = 0  # L105
$ws.Cells.Item(105,13).Value = -49190.25  # M105
$ws.Cells.Item(105,14).ClearContents()  # N105

$ws.Cells.Item(107,8).Value = 2315.3572  # H107
$ws.Cells.Item(107,9).Value = 2176.6086  # I107
$ws.Cells.Item(107,11).Value = 2176.6086  # K107
$ws.Cells.Item(107,13).Value = -256.6086  # M107

$ws.Cells.Item(137,8).Value = 123109.78  # H137
$ws.Cells.Item(137,9).Value = 105797.6  # I137
$ws.Cells.Item(137,11).Value = 105797.6  # K137
$ws.Cells.Item(137,13).Value = -100697.6  # M137

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58,8).Value = 8873184  # H58
$ws.Cells.Item(58,9).Value = 1816.1428  # I58
$ws.Cells.Item(58,10).Value = 16179017  # J58
$ws.Cells.Item(58,11).Value = 1816.1428  # K58
$ws.Cells.Item(58,12).Value = 16179017  # L58
$ws.Cells.Item(58,13).Value = -1613.1428  # M58
$ws.Cells.Item(58,14).Value = -16179423  # N58

$ws.Cells.Item(105,8).Value = 9478.462  # H105
$ws.Cells.Item(105,9).Value = 9560.833000000001  # I105
$ws.Cells.Item(105,11).Value = 9560.833000000001  # K105
$ws.Cells.Item(105,13).Value = -7813.833000000001  # M105

$ws.Cells.Item(132,8).Value = 3518.9023  # H132
$ws.Cells.Item(132,9).Value = 3022.625  # I132
$ws.Cells.Item(132,10).Value = 5283.4443  # J132
$ws.Cells.Item(132,11).Value = 9067.875  # K132
$ws.Cells.Item(132,12).Value = 15850.3329  # L132
$ws.Cells.Item(132,13).Value = -6537.875  # M132
$ws.Cells.Item(132,14).Value = -20910.3329  # N132

$ws.Cells.Item(136,8).Value = 8873184  # H136
$ws.Cells.Item(136,9).Value = 1816.1428  # I136
$ws.Cells.Item(136,10).Value = 16179017  # J136
$ws.Cells.Item(136,11).Value = 5448.428400000001  # K136
$ws.Cells.Item(136,12).Value = 48537051  # L136
$ws.Cells.Item(136,13).Value = -2898.428400000001  # M136
$ws.Cells.Item(136,14).Value = -48542151  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(138,8).Value = 4072.6  # H138
$ws.Cells.Item(138,9).Value = 4594.2856  # I138
$ws.Cells.Item(138,11).Value = 13782.8568  # K138
$ws.Cells.Item(138,13).Value = -8642.856800000001  # M138

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9,8).Value = 463.5  # H9
$ws.Cells.Item(9,9).Value = 326.2  # I9
$ws.Cells.Item(9,10).Value = 1150  # J9
$ws.Cells.Item(9,11).Value = 326.2  # K9
$ws.Cells.Item(9,12).Value = 1150  # L9
$ws.Cells.Item(9,13).Value = -156.2  # M9
$ws.Cells.Item(9,14).Value = -1490  # N9

$ws.Cells.Item(80,8).Value = 2475  # H80
$ws.Cells.Item(80,9).Value = 2200  # I80
$ws.Cells.Item(80,11).Value = 2200  # K80
$ws.Cells.Item(80,13).Value = -1202  # M80

$ws.Cells.Item(83,8).Value = 2475  # H83
$ws.Cells.Item(83,9).Value = 2200  # I83
$ws.Cells.Item(83,11).Value = 11000  # K83
$ws.Cells.Item(83,13).Value = -6008  # M83

$ws.Cells.Item(122,8).Value = 2827  # H122
$ws.Cells.Item(122,9).Value = 2207.1428  # I122
$ws.Cells.Item(122,10).Value = 4996.5  # J122
$ws.Cells.Item(122,11).Value = 6621.428400000001  # K122
$ws.Cells.Item(122,12).Value = 14989.5  # L122
$ws.Cells.Item(122,13).Value = -4171.428400000001  # M122
$ws.Cells.Item(122,14).Value = -19889.5  # N122

$ws.Cells.Item(132,8).Value = 1180.25  # H132
$ws.Cells.Item(132,9).Value = 992.0769  # I132
$ws.Cells.Item(132,11).Value = 2976.2307  # K132
$ws.Cells.Item(132,13).Value = -446.2307000000001  # M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16,8).Value = 2759.75  # H16
$ws.Cells.Item(16,9).Value = 2356.3333  # I16
$ws.Cells.Item(16,11).Value = 2356.3333  # K16
$ws.Cells.Item(16,13).Value = -2186.3333  # M16

$ws.Cells.Item(132,8).Value = 4953.184  # H132
$ws.Cells.Item(132,9).Value = 2586.1333  # I132
$ws.Cells.Item(132,10).Value = 13829.625  # J132
$ws.Cells.Item(132,11).Value = 7758.3999  # K132
$ws.Cells.Item(132,12).Value = 41488.875  # L132
$ws.Cells.Item(132,13).Value = -5228.3999  # M132
$ws.Cells.Item(132,14).Value = -46548.875  # N132

$ws.Cells.Item(133,8).Value = 44175.4  # H133
$ws.Cells.Item(133,10).Value = 59994.5  # J133
$ws.Cells.Item(133,12).Value = 59994.5  # L133
$ws.Cells.Item(133,14).Value = -65054.5  # N133

$ws.Cells.Item(136,8).Value = 2444.0454  # H136
$ws.Cells.Item(136,9).Value = 2313.45  # I136
$ws.Cells.Item(136,11).Value = 6940.349999999999  # K136
$ws.Cells.Item(136,13).Value = -4390.349999999999  # M136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107,8).Value = 903.6  # H107
$ws.Cells.Item(107,9).Value = 684.2105  # I107
$ws.Cells.Item(107,11).Value = 2052.6315  # K107
$ws.Cells.Item(107,13).Value = -132.6315  # M107

$ws.Cells.Item(132,8).Value = 1216.2142  # H132
$ws.Cells.Item(132,9).Value = 852.5  # I132
$ws.Cells.Item(132,11).Value = 2557.5  # K132
$ws.Cells.Item(132,13).Value = -27.5  # M132

$ws.Cells.Item(136,8).Value = 5011.521  # H136
$ws.Cells.Item(136,9).Value = 3029.5134  # I136
$ws.Cells.Item(136,11).Value = 9088.540199999999  # K136
$ws.Cells.Item(136,13).Value = -6538.540199999999  # M136
